# Update the data dictionary export/import sample workbook:
#  - rename the "Data Type" display values to their human-friendly form
#    (date -> Date, plain -> Plain, select -> Multiple Choice)
#  - add a new "other" / "Nonsense" data type row to caseType1 so the
#    fixture exercises an unrecognised data type on import
#  - leave the valid-values sheet untouched content-wise

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # caseType1
$ws2 = $wb.Worksheets.Item(2)   # caseType1-vl
$ws3 = $wb.Worksheets.Item(3)   # caseType2

# Use the "display" values for the existing data types instead of the
# raw/internal ones. caseType1 and caseType2 share the same "Data Type"
# strings (date/plain/select), so both sheets need updating to keep
# them in sync, same as the shared-string table rename.
$ws1.Range("C2").Value = "Date"
$ws1.Range("C3").Value = "Plain"
$ws1.Range("C4").Value = "Multiple Choice"

$ws3.Range("C2").Value = "Date"
$ws3.Range("C3").Value = "Plain"
$ws3.Range("C4").Value = "Multiple Choice"

# Add a new row exercising a bogus/unrecognised data type.
$ws1.Range("A5").Value = "other"
$ws1.Range("C5").Value = "Nonsense"
$ws1.Range("E5").Value = 0

# Restore the selection/active-sheet state: caseType2's cursor moves to
# C5, then caseType1 becomes the active tab with the cursor on the newly
# added row's last cell (E5). Order matters: selecting on caseType2 first
# (making it briefly active), then caseType1 last so it ends up active.
[void]$ws3.Range("C5").Select()
[void]$ws1.Range("E5").Select()
